$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Save" in H1, matching the style of the existing header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill H2:H21 with 0, matching style of existing data column (e.g. G2:G21 - default style)
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
